# Applies the Brynhildr_Profits crafting-leve profit recalculation update.
# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns (H:N) for the affected leve rows across all 8 job sheets.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Cells.Item(32, 8).Value = 1492.0769
$ws.Cells.Item(32, 9).Value = 1774.875
$ws.Cells.Item(32, 10).Value = 1039.6
$ws.Cells.Item(32, 11).Value = 1774.875
$ws.Cells.Item(32, 12).Value = 1039.6
$ws.Cells.Item(32, 13).Value = -1448.875
$ws.Cells.Item(32, 14).Value = -1691.6

# Row 33
$ws.Cells.Item(33, 8).Value = 192.6
$ws.Cells.Item(33, 9).Value = 156.82353
$ws.Cells.Item(33, 11).Value = 156.82353
$ws.Cells.Item(33, 13).Value = 72.17646999999999

# Row 40
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 13).Value = $null

# Row 52
$ws.Cells.Item(52, 8).Value = 4767.091
$ws.Cells.Item(52, 10).Value = 4744.2
$ws.Cells.Item(52, 12).Value = 14232.6
$ws.Cells.Item(52, 14).Value = -14552.6

# Row 64
$ws.Cells.Item(64, 8).Value = 5033.1113
$ws.Cells.Item(64, 9).Value = 3271.6667
$ws.Cells.Item(64, 10).Value = 5913.8335
$ws.Cells.Item(64, 11).Value = 3271.6667
$ws.Cells.Item(64, 12).Value = 5913.8335
$ws.Cells.Item(64, 13).Value = -3023.6667
$ws.Cells.Item(64, 14).Value = -6409.8335

# Row 67
$ws.Cells.Item(67, 8).Value = 5033.1113
$ws.Cells.Item(67, 9).Value = 3271.6667
$ws.Cells.Item(67, 10).Value = 5913.8335
$ws.Cells.Item(67, 11).Value = 3271.6667
$ws.Cells.Item(67, 12).Value = 5913.8335
$ws.Cells.Item(67, 13).Value = -2413.6667
$ws.Cells.Item(67, 14).Value = -7629.8335

# Row 92
$ws.Cells.Item(92, 8).Value = 1041.0416
$ws.Cells.Item(92, 9).Value = 1278
$ws.Cells.Item(92, 10).Value = 465.57144
$ws.Cells.Item(92, 11).Value = 1278
$ws.Cells.Item(92, 12).Value = 465.57144
$ws.Cells.Item(92, 13).Value = -30
$ws.Cells.Item(92, 14).Value = -2961.57144

# Row 106
$ws.Cells.Item(106, 8).Value = 3283.682
$ws.Cells.Item(106, 9).Value = 2912.05
$ws.Cells.Item(106, 11).Value = 2912.05
$ws.Cells.Item(106, 13).Value = -2281.05

# Row 113
$ws.Cells.Item(113, 8).Value = 2633.1333
$ws.Cells.Item(113, 9).Value = 2333.0833
$ws.Cells.Item(113, 11).Value = 2333.0833
$ws.Cells.Item(113, 13).Value = 920.9167000000002

# Row 116
$ws.Cells.Item(116, 8).Value = 12766.125
$ws.Cells.Item(116, 9).Value = 12038.333
$ws.Cells.Item(116, 10).Value = 14949.5
$ws.Cells.Item(116, 11).Value = 12038.333
$ws.Cells.Item(116, 12).Value = 14949.5
$ws.Cells.Item(116, 13).Value = -8596.333000000001
$ws.Cells.Item(116, 14).Value = -21833.5

# Row 132
$ws.Cells.Item(132, 8).Value = 12085.333
$ws.Cells.Item(132, 9).Value = 12589.6
$ws.Cells.Item(132, 11).Value = 37768.8
$ws.Cells.Item(132, 13).Value = -35238.8

# Row 137
$ws.Cells.Item(137, 8).Value = 15635653
$ws.Cells.Item(137, 9).Value = 35717704
$ws.Cells.Item(137, 10).Value = 16279.611
$ws.Cells.Item(137, 11).Value = 107153112
$ws.Cells.Item(137, 12).Value = 48838.833
$ws.Cells.Item(137, 13).Value = -107150562
$ws.Cells.Item(137, 14).Value = -53938.833

# Row 138
$ws.Cells.Item(138, 8).Value = 8202.191999999999
$ws.Cells.Item(138, 9).Value = 9649.75
$ws.Cells.Item(138, 10).Value = 7939
$ws.Cells.Item(138, 11).Value = 28949.25
$ws.Cells.Item(138, 12).Value = 23817
$ws.Cells.Item(138, 13).Value = -23809.25
$ws.Cells.Item(138, 14).Value = -34097

# Row 140
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).Value = $null

# Row 141
$ws.Cells.Item(141, 8).Value = 7221.393
$ws.Cells.Item(141, 9).Value = 2350
$ws.Cells.Item(141, 10).Value = 12092.786
$ws.Cells.Item(141, 11).Value = 7050
$ws.Cells.Item(141, 12).Value = 36278.358
$ws.Cells.Item(141, 13).Value = -1870
$ws.Cells.Item(141, 14).Value = -46638.358

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 1208
$ws.Cells.Item(2, 9).Value = 1166.25
$ws.Cells.Item(2, 10).Value = 1308.2
$ws.Cells.Item(2, 11).Value = 1166.25
$ws.Cells.Item(2, 12).Value = 1308.2
$ws.Cells.Item(2, 13).Value = -1053.25
$ws.Cells.Item(2, 14).Value = -1534.2

# Row 4
$ws.Cells.Item(4, 8).Value = 300
$ws.Cells.Item(4, 9).Value = 300
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 300
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = -184
$ws.Cells.Item(4, 14).Value = $null

# Row 55
$ws.Cells.Item(55, 8).Value = 17620
$ws.Cells.Item(55, 10).Value = 28700
$ws.Cells.Item(55, 12).Value = 28700
$ws.Cells.Item(55, 14).Value = -29330

# Row 61
$ws.Cells.Item(61, 8).Value = 2328767.2
$ws.Cells.Item(61, 9).Value = 2927.2856
$ws.Cells.Item(61, 10).Value = 6670335.5
$ws.Cells.Item(61, 11).Value = 2927.2856
$ws.Cells.Item(61, 12).Value = 6670335.5
$ws.Cells.Item(61, 13).Value = -2715.2856
$ws.Cells.Item(61, 14).Value = -6670759.5

# Row 110
$ws.Cells.Item(110, 8).Value = 1109.8182
$ws.Cells.Item(110, 9).Value = 782.4
$ws.Cells.Item(110, 10).Value = 1811.4286
$ws.Cells.Item(110, 11).Value = 782.4
$ws.Cells.Item(110, 12).Value = 1811.4286
$ws.Cells.Item(110, 13).Value = 1262.6
$ws.Cells.Item(110, 14).Value = -5901.4286

# Row 116
$ws.Cells.Item(116, 8).Value = 1208
$ws.Cells.Item(116, 9).Value = 1166.25
$ws.Cells.Item(116, 10).Value = 1308.2
$ws.Cells.Item(116, 11).Value = 1166.25
$ws.Cells.Item(116, 12).Value = 1308.2
$ws.Cells.Item(116, 13).Value = 1127.75
$ws.Cells.Item(116, 14).Value = -5896.2

# Row 122
$ws.Cells.Item(122, 8).Value = 1275.3704
$ws.Cells.Item(122, 9).Value = 1093.909
$ws.Cells.Item(122, 11).Value = 3281.727
$ws.Cells.Item(122, 13).Value = -831.7270000000003

# Row 132
$ws.Cells.Item(132, 8).Value = 4299.7666
$ws.Cells.Item(132, 9).Value = 2916.6667
$ws.Cells.Item(132, 10).Value = 5990.222
$ws.Cells.Item(132, 11).Value = 8750.000100000001
$ws.Cells.Item(132, 12).Value = 17970.666
$ws.Cells.Item(132, 13).Value = -6220.000100000001
$ws.Cells.Item(132, 14).Value = -23030.666

# Row 136
$ws.Cells.Item(136, 8).Value = 2328767.2
$ws.Cells.Item(136, 9).Value = 2927.2856
$ws.Cells.Item(136, 10).Value = 6670335.5
$ws.Cells.Item(136, 11).Value = 8781.856800000001
$ws.Cells.Item(136, 12).Value = 20011006.5
$ws.Cells.Item(136, 13).Value = -6231.856800000001
$ws.Cells.Item(136, 14).Value = -20016106.5

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 1208
$ws.Cells.Item(3, 9).Value = 1166.25
$ws.Cells.Item(3, 10).Value = 1308.2
$ws.Cells.Item(3, 11).Value = 1166.25
$ws.Cells.Item(3, 12).Value = 1308.2
$ws.Cells.Item(3, 13).Value = -1052.25
$ws.Cells.Item(3, 14).Value = -1536.2

# Row 20
$ws.Cells.Item(20, 8).Value = 35150.59
$ws.Cells.Item(20, 9).Value = 41751.832
$ws.Cells.Item(20, 10).Value = 13146.444
$ws.Cells.Item(20, 11).Value = 41751.832
$ws.Cells.Item(20, 12).Value = 13146.444
$ws.Cells.Item(20, 13).Value = -41504.832
$ws.Cells.Item(20, 14).Value = -13640.444

# Row 22
$ws.Cells.Item(22, 8).Value = 298.5
$ws.Cells.Item(22, 9).Value = 298.2
$ws.Cells.Item(22, 10).Value = 300
$ws.Cells.Item(22, 11).Value = 298.2
$ws.Cells.Item(22, 12).Value = 300
$ws.Cells.Item(22, 13).Value = -125.2
$ws.Cells.Item(22, 14).Value = -646

# Row 94
$ws.Cells.Item(94, 8).Value = 3294.75
$ws.Cells.Item(94, 9).Value = 3276.6843
$ws.Cells.Item(94, 10).Value = 3363.4
$ws.Cells.Item(94, 11).Value = 3276.6843
$ws.Cells.Item(94, 12).Value = 3363.4
$ws.Cells.Item(94, 13).Value = -2825.6843
$ws.Cells.Item(94, 14).Value = -4265.4

# Row 99
$ws.Cells.Item(99, 8).Value = 5416.5
$ws.Cells.Item(99, 9).Value = 6062.8335
$ws.Cells.Item(99, 10).Value = 2831.1667
$ws.Cells.Item(99, 11).Value = 6062.8335
$ws.Cells.Item(99, 12).Value = 2831.1667
$ws.Cells.Item(99, 13).Value = -4564.8335
$ws.Cells.Item(99, 14).Value = -5827.1667

# Row 105
$ws.Cells.Item(105, 8).Value = 1992.0714
$ws.Cells.Item(105, 9).Value = 1847.9565
$ws.Cells.Item(105, 10).Value = 2655
$ws.Cells.Item(105, 11).Value = 1847.9565
$ws.Cells.Item(105, 12).Value = 2655
$ws.Cells.Item(105, 13).Value = -100.9565
$ws.Cells.Item(105, 14).Value = -6149

# Row 107
$ws.Cells.Item(107, 8).Value = 1767.8948
$ws.Cells.Item(107, 9).Value = 1522.875
$ws.Cells.Item(107, 10).Value = 1946.091
$ws.Cells.Item(107, 11).Value = 1522.875
$ws.Cells.Item(107, 12).Value = 1946.091
$ws.Cells.Item(107, 13).Value = 397.125
$ws.Cells.Item(107, 14).Value = -5786.091

# Row 134
$ws.Cells.Item(134, 8).Value = 2781336.8
$ws.Cells.Item(134, 9).Value = 3604.45
$ws.Cells.Item(134, 10).Value = 8336801.5
$ws.Cells.Item(134, 11).Value = 10813.35
$ws.Cells.Item(134, 12).Value = 25010404.5
$ws.Cells.Item(134, 13).Value = -8278.349999999999
$ws.Cells.Item(134, 14).Value = -25015474.5

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 28036.5
$ws.Cells.Item(16, 9).Value = 33914.5
$ws.Cells.Item(16, 10).Value = 4524.5
$ws.Cells.Item(16, 11).Value = 33914.5
$ws.Cells.Item(16, 12).Value = 4524.5
$ws.Cells.Item(16, 13).Value = -33627.5
$ws.Cells.Item(16, 14).Value = -5098.5

# Row 22
$ws.Cells.Item(22, 8).Value = 341.1111
$ws.Cells.Item(22, 9).Value = 209.23077
$ws.Cells.Item(22, 10).Value = 684
$ws.Cells.Item(22, 11).Value = 209.23077
$ws.Cells.Item(22, 12).Value = 684
$ws.Cells.Item(22, 13).Value = 140.76923
$ws.Cells.Item(22, 14).Value = -1384

# Row 31
$ws.Cells.Item(31, 8).Value = 567015.5600000001
$ws.Cells.Item(31, 9).Value = 754423.75
$ws.Cells.Item(31, 10).Value = 4791.0835
$ws.Cells.Item(31, 11).Value = 754423.75
$ws.Cells.Item(31, 12).Value = 4791.0835
$ws.Cells.Item(31, 13).Value = -754128.75
$ws.Cells.Item(31, 14).Value = -5381.0835

# Row 34
$ws.Cells.Item(34, 8).Value = 567015.5600000001
$ws.Cells.Item(34, 9).Value = 754423.75
$ws.Cells.Item(34, 10).Value = 4791.0835
$ws.Cells.Item(34, 11).Value = 754423.75
$ws.Cells.Item(34, 12).Value = 4791.0835
$ws.Cells.Item(34, 13).Value = -754221.75
$ws.Cells.Item(34, 14).Value = -5195.0835

# Row 50
$ws.Cells.Item(50, 8).Value = 27533.334
$ws.Cells.Item(50, 9).Value = 25000
$ws.Cells.Item(50, 10).Value = 28800
$ws.Cells.Item(50, 11).Value = 25000
$ws.Cells.Item(50, 12).Value = 28800
$ws.Cells.Item(50, 13).Value = -24375
$ws.Cells.Item(50, 14).Value = -30050

# Row 51
$ws.Cells.Item(51, 8).Value = 28800
$ws.Cells.Item(51, 10).Value = 28800
$ws.Cells.Item(51, 12).Value = 28800
$ws.Cells.Item(51, 14).Value = -30272

# Row 58
$ws.Cells.Item(58, 8).Value = 5128739
$ws.Cells.Item(58, 9).Value = 12823776
$ws.Cells.Item(58, 10).Value = 1556043
$ws.Cells.Item(58, 11).Value = 12823776
$ws.Cells.Item(58, 12).Value = 1556043
$ws.Cells.Item(58, 13).Value = -12823573
$ws.Cells.Item(58, 14).Value = -1556449

# Row 59
$ws.Cells.Item(59, 8).Value = 38108.8
$ws.Cells.Item(59, 9).Value = 32497.5
$ws.Cells.Item(59, 10).Value = 41849.668
$ws.Cells.Item(59, 11).Value = 32497.5
$ws.Cells.Item(59, 12).Value = 41849.668
$ws.Cells.Item(59, 13).Value = -31352.5
$ws.Cells.Item(59, 14).Value = -44139.668

# Row 60
$ws.Cells.Item(60, 8).Value = 28799.666
$ws.Cells.Item(60, 10).Value = 28799.666
$ws.Cells.Item(60, 12).Value = 28799.666
$ws.Cells.Item(60, 14).Value = -29821.666

# Row 61
$ws.Cells.Item(61, 8).Value = 28800
$ws.Cells.Item(61, 10).Value = 28800
$ws.Cells.Item(61, 12).Value = 28800
$ws.Cells.Item(61, 14).Value = -29496

# Row 68
$ws.Cells.Item(68, 8).Value = 43500
$ws.Cells.Item(68, 10).Value = 43500
$ws.Cells.Item(68, 12).Value = 43500
$ws.Cells.Item(68, 14).Value = -44998

# Row 69
$ws.Cells.Item(69, 8).Value = 11250
$ws.Cells.Item(69, 9).Value = 7500
$ws.Cells.Item(69, 11).Value = 7500
$ws.Cells.Item(69, 13).Value = -6751

# Row 71
$ws.Cells.Item(71, 8).Value = 43500
$ws.Cells.Item(71, 10).Value = 43500
$ws.Cells.Item(71, 12).Value = 130500
$ws.Cells.Item(71, 14).Value = -137988

# Row 72
$ws.Cells.Item(72, 8).Value = 11250
$ws.Cells.Item(72, 9).Value = 7500
$ws.Cells.Item(72, 11).Value = 22500
$ws.Cells.Item(72, 13).Value = -18756

# Row 74
$ws.Cells.Item(74, 8).Value = 43642.5
$ws.Cells.Item(74, 10).Value = 45000
$ws.Cells.Item(74, 12).Value = 45000
$ws.Cells.Item(74, 14).Value = -46748

# Row 77
$ws.Cells.Item(77, 8).Value = 43642.5
$ws.Cells.Item(77, 10).Value = 45000
$ws.Cells.Item(77, 12).Value = 135000
$ws.Cells.Item(77, 14).Value = -143736

# Row 93
$ws.Cells.Item(93, 8).Value = 17599.5
$ws.Cells.Item(93, 9).Value = 17599.5
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 17599.5
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = -15727.5
$ws.Cells.Item(93, 14).Value = $null

# Row 105
$ws.Cells.Item(105, 8).Value = 16140.85
$ws.Cells.Item(105, 9).Value = 16140.85
$ws.Cells.Item(105, 11).Value = 16140.85
$ws.Cells.Item(105, 13).Value = -14393.85

# Row 107
$ws.Cells.Item(107, 8).Value = 534.4
$ws.Cells.Item(107, 9).Value = 538.6923
$ws.Cells.Item(107, 10).Value = 506.5
$ws.Cells.Item(107, 11).Value = 538.6923
$ws.Cells.Item(107, 12).Value = 506.5
$ws.Cells.Item(107, 13).Value = 1381.3077
$ws.Cells.Item(107, 14).Value = -4346.5

# Row 113
$ws.Cells.Item(113, 8).Value = 28036.5
$ws.Cells.Item(113, 9).Value = 33914.5
$ws.Cells.Item(113, 10).Value = 4524.5
$ws.Cells.Item(113, 11).Value = 33914.5
$ws.Cells.Item(113, 12).Value = 4524.5
$ws.Cells.Item(113, 13).Value = -31744.5
$ws.Cells.Item(113, 14).Value = -8864.5

# Row 131
$ws.Cells.Item(131, 8).Value = 59839
$ws.Cells.Item(131, 10).Value = 59839
$ws.Cells.Item(131, 12).Value = 59839
$ws.Cells.Item(131, 14).Value = -69919

# Row 132
$ws.Cells.Item(132, 8).Value = 3448.125
$ws.Cells.Item(132, 9).Value = 3336.5386
$ws.Cells.Item(132, 11).Value = 10009.6158
$ws.Cells.Item(132, 13).Value = -7479.6158

# Row 136
$ws.Cells.Item(136, 8).Value = 5128739
$ws.Cells.Item(136, 9).Value = 12823776
$ws.Cells.Item(136, 10).Value = 1556043
$ws.Cells.Item(136, 11).Value = 38471328
$ws.Cells.Item(136, 12).Value = 4668129
$ws.Cells.Item(136, 13).Value = -38468778
$ws.Cells.Item(136, 14).Value = -4673229

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Cells.Item(7, 8).Value = 126.72727
$ws.Cells.Item(7, 9).Value = 29.666666
$ws.Cells.Item(7, 10).Value = 243.2
$ws.Cells.Item(7, 11).Value = 88.99999800000001
$ws.Cells.Item(7, 12).Value = 729.5999999999999
$ws.Cells.Item(7, 13).Value = 23.00000199999999
$ws.Cells.Item(7, 14).Value = -953.5999999999999

# Row 32
$ws.Cells.Item(32, 8).Value = 2503699.8
$ws.Cells.Item(32, 9).Value = 4999
$ws.Cells.Item(32, 10).Value = 3336600
$ws.Cells.Item(32, 11).Value = 14997
$ws.Cells.Item(32, 12).Value = 10009800
$ws.Cells.Item(32, 13).Value = -14714
$ws.Cells.Item(32, 14).Value = -10010366

# Row 46
$ws.Cells.Item(46, 8).Value = 834513
$ws.Cells.Item(46, 9).Value = 930.8
$ws.Cells.Item(46, 10).Value = 1429928.9
$ws.Cells.Item(46, 11).Value = 2792.4
$ws.Cells.Item(46, 12).Value = 4289786.699999999
$ws.Cells.Item(46, 13).Value = -2701.4
$ws.Cells.Item(46, 14).Value = -4289968.699999999

# Row 58
$ws.Cells.Item(58, 8).Value = 9901
$ws.Cells.Item(58, 9).Value = 7002.5
$ws.Cells.Item(58, 10).Value = 11833.333
$ws.Cells.Item(58, 11).Value = 21007.5
$ws.Cells.Item(58, 12).Value = 35499.999
$ws.Cells.Item(58, 13).Value = -20879.5
$ws.Cells.Item(58, 14).Value = -35755.999

# Row 107
$ws.Cells.Item(107, 8).Value = 6611.0713
$ws.Cells.Item(107, 10).Value = 7883.5557
$ws.Cells.Item(107, 12).Value = 23650.6671
$ws.Cells.Item(107, 14).Value = -27490.6671

# Row 138
$ws.Cells.Item(138, 8).Value = 16169.568
$ws.Cells.Item(138, 9).Value = 19761.834
$ws.Cells.Item(138, 10).Value = 13682.615
$ws.Cells.Item(138, 11).Value = 59285.50199999999
$ws.Cells.Item(138, 12).Value = 41047.845
$ws.Cells.Item(138, 13).Value = -54145.50199999999
$ws.Cells.Item(138, 14).Value = -51327.845

# Row 139
$ws.Cells.Item(139, 8).Value = 5614.7334
$ws.Cells.Item(139, 9).Value = 3140.0952
$ws.Cells.Item(139, 10).Value = 11388.889
$ws.Cells.Item(139, 11).Value = 9420.285600000001
$ws.Cells.Item(139, 12).Value = 34166.667
$ws.Cells.Item(139, 13).Value = -4280.285600000001
$ws.Cells.Item(139, 14).Value = -44446.667

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Cells.Item(43, 8).Value = 2517
$ws.Cells.Item(43, 9).Value = 2517
$ws.Cells.Item(43, 11).Value = 2517
$ws.Cells.Item(43, 13).Value = -2366

# Row 70
$ws.Cells.Item(70, 8).Value = 16199.046
$ws.Cells.Item(70, 9).Value = 25375.666
$ws.Cells.Item(70, 10).Value = 9846
$ws.Cells.Item(70, 11).Value = 25375.666
$ws.Cells.Item(70, 12).Value = 9846
$ws.Cells.Item(70, 13).Value = -25105.666
$ws.Cells.Item(70, 14).Value = -10386

# Row 73
$ws.Cells.Item(73, 8).Value = 16199.046
$ws.Cells.Item(73, 9).Value = 25375.666
$ws.Cells.Item(73, 10).Value = 9846
$ws.Cells.Item(73, 11).Value = 25375.666
$ws.Cells.Item(73, 12).Value = 9846
$ws.Cells.Item(73, 13).Value = -24439.666
$ws.Cells.Item(73, 14).Value = -11718

# Row 97
$ws.Cells.Item(97, 8).Value = 966.0833
$ws.Cells.Item(97, 9).Value = 841.5909
$ws.Cells.Item(97, 10).Value = 1161.7142
$ws.Cells.Item(97, 11).Value = 841.5909
$ws.Cells.Item(97, 12).Value = 1161.7142
$ws.Cells.Item(97, 13).Value = -345.5909
$ws.Cells.Item(97, 14).Value = -2153.7142

# Row 122
$ws.Cells.Item(122, 8).Value = 73380.64
$ws.Cells.Item(122, 9).Value = 84860.75
$ws.Cells.Item(122, 11).Value = 254582.25
$ws.Cells.Item(122, 13).Value = -252132.25

# Row 132
$ws.Cells.Item(132, 8).Value = 8147.035
$ws.Cells.Item(132, 9).Value = 9440.679
$ws.Cells.Item(132, 10).Value = 6898
$ws.Cells.Item(132, 11).Value = 28322.037
$ws.Cells.Item(132, 12).Value = 20694
$ws.Cells.Item(132, 13).Value = -25792.037
$ws.Cells.Item(132, 14).Value = -25754

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 1204.6364
$ws.Cells.Item(22, 9).Value = 833.3333
$ws.Cells.Item(22, 10).Value = 1343.875
$ws.Cells.Item(22, 11).Value = 833.3333
$ws.Cells.Item(22, 12).Value = 1343.875
$ws.Cells.Item(22, 13).Value = -538.3333
$ws.Cells.Item(22, 14).Value = -1933.875

# Row 27
$ws.Cells.Item(27, 8).Value = 1204.6364
$ws.Cells.Item(27, 9).Value = 833.3333
$ws.Cells.Item(27, 10).Value = 1343.875
$ws.Cells.Item(27, 11).Value = 833.3333
$ws.Cells.Item(27, 12).Value = 1343.875
$ws.Cells.Item(27, 13).Value = -726.3333
$ws.Cells.Item(27, 14).Value = -1557.875

# Row 46
$ws.Cells.Item(46, 8).Value = 3999
$ws.Cells.Item(46, 9).Value = 995
$ws.Cells.Item(46, 10).Value = 4750
$ws.Cells.Item(46, 11).Value = 995
$ws.Cells.Item(46, 12).Value = 4750
$ws.Cells.Item(46, 13).Value = -807
$ws.Cells.Item(46, 14).Value = -5126

# Row 47
$ws.Cells.Item(47, 8).Value = 26997
$ws.Cells.Item(47, 10).Value = 26997
$ws.Cells.Item(47, 12).Value = 26997
$ws.Cells.Item(47, 14).Value = -27977

# Row 52
$ws.Cells.Item(52, 8).Value = 26997
$ws.Cells.Item(52, 10).Value = 26997
$ws.Cells.Item(52, 12).Value = 26997
$ws.Cells.Item(52, 14).Value = -27463

# Row 55
$ws.Cells.Item(55, 8).Value = 363.12122
$ws.Cells.Item(55, 10).Value = 599.0625
$ws.Cells.Item(55, 12).Value = 599.0625
$ws.Cells.Item(55, 14).Value = -945.0625

# Row 61
$ws.Cells.Item(61, 8).Value = 12817.272
$ws.Cells.Item(61, 9).Value = 15571.857
$ws.Cells.Item(61, 11).Value = 15571.857
$ws.Cells.Item(61, 13).Value = -15369.857

# Row 93
$ws.Cells.Item(93, 8).Value = 2362.5625
$ws.Cells.Item(93, 9).Value = 1274.9412
$ws.Cells.Item(93, 11).Value = 1274.9412
$ws.Cells.Item(93, 13).Value = -26.94119999999998

# Row 113
$ws.Cells.Item(113, 8).Value = 12817.272
$ws.Cells.Item(113, 9).Value = 15571.857
$ws.Cells.Item(113, 11).Value = 15571.857
$ws.Cells.Item(113, 13).Value = -13401.857

# Row 114
$ws.Cells.Item(114, 8).Value = 40000
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 14).Value = $null

# Row 132
$ws.Cells.Item(132, 8).Value = 4389786
$ws.Cells.Item(132, 9).Value = 4633390.5
$ws.Cells.Item(132, 10).Value = 4905
$ws.Cells.Item(132, 11).Value = 13900171.5
$ws.Cells.Item(132, 12).Value = 14715
$ws.Cells.Item(132, 13).Value = -13897641.5
$ws.Cells.Item(132, 14).Value = -19775

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Cells.Item(45, 8).Value = 40882
$ws.Cells.Item(45, 10).Value = 40196.668
$ws.Cells.Item(45, 12).Value = 40196.668
$ws.Cells.Item(45, 14).Value = -41178.668

# Row 62
$ws.Cells.Item(62, 8).Value = 29165.334
$ws.Cells.Item(62, 9).Value = 17997.5
$ws.Cells.Item(62, 10).Value = 51501
$ws.Cells.Item(62, 11).Value = 17997.5
$ws.Cells.Item(62, 12).Value = 51501
$ws.Cells.Item(62, 13).Value = -17373.5
$ws.Cells.Item(62, 14).Value = -52749

# Row 65
$ws.Cells.Item(65, 8).Value = 29165.334
$ws.Cells.Item(65, 9).Value = 17997.5
$ws.Cells.Item(65, 10).Value = 51501
$ws.Cells.Item(65, 11).Value = 89987.5
$ws.Cells.Item(65, 12).Value = 257505
$ws.Cells.Item(65, 13).Value = -86867.5
$ws.Cells.Item(65, 14).Value = -263745

# Row 69
$ws.Cells.Item(69, 8).Value = 52387.5
$ws.Cells.Item(69, 9).Value = 35000
$ws.Cells.Item(69, 10).Value = 58183.332
$ws.Cells.Item(69, 11).Value = 35000
$ws.Cells.Item(69, 12).Value = 58183.332
$ws.Cells.Item(69, 13).Value = -34251
$ws.Cells.Item(69, 14).Value = -59681.332

# Row 72
$ws.Cells.Item(72, 8).Value = 52387.5
$ws.Cells.Item(72, 9).Value = 35000
$ws.Cells.Item(72, 10).Value = 58183.332
$ws.Cells.Item(72, 11).Value = 105000
$ws.Cells.Item(72, 12).Value = 174549.996
$ws.Cells.Item(72, 13).Value = -101256
$ws.Cells.Item(72, 14).Value = -182037.996

# Row 132
$ws.Cells.Item(132, 8).Value = 3970278.5
$ws.Cells.Item(132, 9).Value = 4275479.5
$ws.Cells.Item(132, 11).Value = 12826438.5
$ws.Cells.Item(132, 13).Value = -12823908.5

# Row 136
$ws.Cells.Item(136, 8).Value = 2257243.2
$ws.Cells.Item(136, 9).Value = 1673713.5
$ws.Cells.Item(136, 10).Value = 2979708.8
$ws.Cells.Item(136, 11).Value = 5021140.5
$ws.Cells.Item(136, 12).Value = 8939126.399999999
$ws.Cells.Item(136, 13).Value = -5018590.5
$ws.Cells.Item(136, 14).Value = -8944226.399999999

Write-Output "Brynhildr_Profits: updated $(86) leve-profit rows across 8 sheets"